$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (Changed) date column C for all existing data rows
#    (rows 2 through 514) from 45202 to 45203.
$ws.Range("C2:C514").Value = 45203

# 2) Give row 514 an explicit row height (matches new rows being appended below it).
$ws.Rows.Item(514).RowHeight = 15

# 3) Append three new data rows (515, 516, 517) for new logging notifications.

# --- Row 515 ---
$ws.Range("A515").Value = "A 47422-2023"
$ws.Range("B515").Value = 45202
$ws.Range("C515").Value = 45203
$ws.Range("D515").Value = "HALLANDS LÄN"
$ws.Range("E515").Value = "FALKENBERG"
$ws.Range("G515").Value = 1
$ws.Range("H515").Value = 0
$ws.Range("I515").Value = 0
$ws.Range("J515").Value = 0
$ws.Range("K515").Value = 0
$ws.Range("L515").Value = 0
$ws.Range("M515").Value = 0
$ws.Range("N515").Value = 0
$ws.Range("O515").Value = 0
$ws.Range("P515").Value = 0
$ws.Range("Q515").Value = 0
$ws.Range("B515:C515").NumberFormat = "YYYY-MM-DD"
$ws.Range("R515").WrapText = $true
$ws.Rows.Item(515).RowHeight = 15

# --- Row 516 ---
$ws.Range("A516").Value = "A 47423-2023"
$ws.Range("B516").Value = 45202
$ws.Range("C516").Value = 45203
$ws.Range("D516").Value = "HALLANDS LÄN"
$ws.Range("E516").Value = "FALKENBERG"
$ws.Range("G516").Value = 0.5
$ws.Range("H516").Value = 0
$ws.Range("I516").Value = 0
$ws.Range("J516").Value = 0
$ws.Range("K516").Value = 0
$ws.Range("L516").Value = 0
$ws.Range("M516").Value = 0
$ws.Range("N516").Value = 0
$ws.Range("O516").Value = 0
$ws.Range("P516").Value = 0
$ws.Range("Q516").Value = 0
$ws.Range("B516:C516").NumberFormat = "YYYY-MM-DD"
$ws.Range("R516").WrapText = $true
$ws.Rows.Item(516).RowHeight = 15

# --- Row 517 ---
$ws.Range("A517").Value = "A 47412-2023"
$ws.Range("B517").Value = 45202
$ws.Range("C517").Value = 45203
$ws.Range("D517").Value = "HALLANDS LÄN"
$ws.Range("E517").Value = "FALKENBERG"
$ws.Range("G517").Value = 2.6
$ws.Range("H517").Value = 0
$ws.Range("I517").Value = 0
$ws.Range("J517").Value = 0
$ws.Range("K517").Value = 0
$ws.Range("L517").Value = 0
$ws.Range("M517").Value = 0
$ws.Range("N517").Value = 0
$ws.Range("O517").Value = 0
$ws.Range("P517").Value = 0
$ws.Range("Q517").Value = 0
$ws.Range("B517:C517").NumberFormat = "YYYY-MM-DD"
$ws.Range("R517").WrapText = $true
